$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (number format, font, borders, alignment) from the last
# existing data row (69) down into the two new rows (70, 71) before writing
# values, so the new rows match the existing style (bold/centered index
# column, date-time number format on the match-date column, etc.). Only the
# used columns (A:V) are copied so the sheet's used range / dimension isn't
# blown out to the full row width.
$ws.Range("A69:V69").Copy()
$ws.Range("A70:V70").PasteSpecial(-4122)
$ws.Range("A69:V69").Copy()
$ws.Range("A71:V71").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Row 70 ----
$ws.Range("A70").Value = 69
$ws.Range("B70").Value = "croatia"
$ws.Range("C70").Value = "hnl"
$ws.Range("D70").Value = "2023-2024"
$ws.Range("E70").Value = 45241.66666666666
$ws.Range("F70").Value = "Rudes"
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = "Istra 1961"
$ws.Range("I70").Value = 4
$ws.Range("J70").Value = 2.62
$ws.Range("K70").Value = "05/11/2023 17:12"
$ws.Range("L70").Value = 2.59
$ws.Range("M70").Value = "11/11/2023 15:59"
$ws.Range("N70").Value = 3.22
$ws.Range("O70").Value = "05/11/2023 17:12"
$ws.Range("P70").Value = 3.04
$ws.Range("Q70").Value = "11/11/2023 15:59"
$ws.Range("R70").Value = 2.78
$ws.Range("S70").Value = "05/11/2023 17:12"
$ws.Range("T70").Value = 3.04
$ws.Range("U70").Value = "11/11/2023 15:59"
$ws.Range("V70").Value = "https://www.betexplorer.com/football/croatia/hnl/rudes-istra-1961/tjhwUcJI/"

# ---- Row 71 ----
$ws.Range("A71").Value = 70
$ws.Range("B71").Value = "croatia"
$ws.Range("C71").Value = "hnl"
$ws.Range("D71").Value = "2023-2024"
$ws.Range("E71").Value = 45241.75694444445
$ws.Range("F71").Value = "Osijek"
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = "Lok. Zagreb"
$ws.Range("I71").Value = 1
$ws.Range("J71").Value = 1.75
$ws.Range("K71").Value = "08/11/2023 17:13"
$ws.Range("L71").Value = 1.7
$ws.Range("M71").Value = "11/11/2023 18:06"
$ws.Range("N71").Value = 3.78
$ws.Range("O71").Value = "08/11/2023 17:13"
$ws.Range("P71").Value = 3.82
$ws.Range("Q71").Value = "11/11/2023 18:06"
$ws.Range("R71").Value = 4.17
$ws.Range("S71").Value = "08/11/2023 17:13"
$ws.Range("T71").Value = 5.03
$ws.Range("U71").Value = "11/11/2023 18:02"
$ws.Range("V71").Value = "https://www.betexplorer.com/football/croatia/hnl/osijek-lok-zagreb/KKgkyFQo/"
